$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

function Find-ShapeById($shapes, [int]$id) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $candidate = $shapes.Item($i)
        if ($candidate.Id -eq $id) {
            return $candidate
        }
    }
    return $null
}

function Replace-Substring($shape, [string]$old, [string]$new) {
    $tr = $shape.TextFrame.TextRange
    $full = $tr.Text
    $idx = $full.IndexOf($old)
    if ($idx -lt 0) {
        Write-Host ("NOT FOUND: '" + $old + "' in shape " + $shape.Id)
        return
    }
    $sub = $tr.Characters($idx + 1, $old.Length)
    $sub.Text = $new
}

# Rectangle 8 (id 46): "VersionedAddressBook" -> "VersionedClinicIo"
$shpVersioned = Find-ShapeById $s.Shapes 46
Replace-Substring $shpVersioned "VersionedAddressBook" "VersionedClinicIo"

# Rectangle 8 (id 92): "<<interface>>ReadOnlyAddressBook" -> "<<interface>>ReadOnlyClinicIo"
$shpReadOnly = Find-ShapeById $s.Shapes 92
Replace-Substring $shpReadOnly "ReadOnlyAddressBook" "ReadOnlyClinicIo"

# Rectangle 8 (id 93): "AddressBook" -> "ClinicIo"
$shpPlain = Find-ShapeById $s.Shapes 93
Replace-Substring $shpPlain "AddressBook" "ClinicIo"
